$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values so Excel does not
# silently convert them to real numbers (matches original inlineStr text cells).
$numericLooking = @(
    "D4",
    "D6",
    "D7",
    "D8",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D17",
    "D18",
    "D19",
    "D22",
    "D23",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value2 = '22.377.68'
$ws.Range("E2").Value2 = '  -0.11%  '

$ws.Range("D3").Value2 = '1.570.29'
$ws.Range("E3").Value2 = '  -0.22%  '

$ws.Range("D4").Value2 = '1.003'
$ws.Range("E4").Value2 = '  +0.20%  '

$ws.Range("E5").Value2 = '  +0.13%  '

$ws.Range("D6").Value2 = '291.20'
$ws.Range("E6").Value2 = '  +0.50%  '

$ws.Range("D7").Value2 = '0.3763'
$ws.Range("E7").Value2 = '  +2.40%  '

$ws.Range("D8").Value2 = '49.68'
$ws.Range("E8").Value2 = '  +0.64%  '

$ws.Range("E9").Value2 = '  +0.91%  '

$ws.Range("D10").Value2 = '0.07613'
$ws.Range("E10").Value2 = '  -0.12%  '

$ws.Range("D11").Value2 = '1.141'
$ws.Range("E11").Value2 = '  -2.47%  '

$ws.Range("E12").Value2 = '  +0.20%  '

$ws.Range("D13").Value2 = '21.14'
$ws.Range("E13").Value2 = '  -0.94%  '

$ws.Range("D14").Value2 = '5.992'
$ws.Range("E14").Value2 = '  -1.20%  '

$ws.Range("D15").Value2 = '6.929'
$ws.Range("E15").Value2 = '  -0.12%  '

$ws.Range("D16").Value2 = '1.570.79'
$ws.Range("E16").Value2 = '  +0.46%  '

$ws.Range("D17").Value2 = '0.00001133'

$ws.Range("D18").Value2 = '90.26'
$ws.Range("E18").Value2 = '  +1.17%  '

$ws.Range("D19").Value2 = '0.06748'
$ws.Range("E19").Value2 = '  +0.01%  '

$ws.Range("E20").Value2 = '  +0.10%  '

$ws.Range("E21").Value2 = '  +0.82%  '

$ws.Range("D22").Value2 = '6.189'
$ws.Range("E22").Value2 = '  -0.97%  '

$ws.Range("D23").Value2 = '11.99'
$ws.Range("E23").Value2 = '  -0.03%  '

$ws.Range("D24").Value2 = '22.384.02'

$ws.Range("E25").Value2 = '  +0.57%  '

$ws.Range("D26").Value2 = '2.650'
$ws.Range("E26").Value2 = '  -10.88%  '

$ws.Range("D27").Value2 = '20.10'
$ws.Range("E27").Value2 = '  +0.44%  '

$ws.Range("D28").Value2 = '147.11'
$ws.Range("E28").Value2 = '  +0.84%  '

$ws.Range("D29").Value2 = '5.044'
$ws.Range("E29").Value2 = '  +1.29%  '

$ws.Range("D30").Value2 = '126.61'

$ws.Range("D31").Value2 = '1.746.33'
$ws.Range("E31").Value2 = '  +0.13%  '

$ws.Range("D32").Value2 = '2.011'
$ws.Range("E32").Value2 = '  +0.33%  '

$ws.Range("D33").Value2 = '6.082'
$ws.Range("E33").Value2 = '  -3.57%  '

$ws.Range("D34").Value2 = '0.9890'
$ws.Range("E34").Value2 = '  -5.40%  '

$ws.Range("D35").Value2 = '10.13'
$ws.Range("E35").Value2 = '  -1.91%  '

$ws.Range("D36").Value2 = '0.08514'
$ws.Range("E36").Value2 = '  +0.74%  '

$ws.Range("D37").Value2 = '0.02530'
$ws.Range("E37").Value2 = '  -0.44%  '

$ws.Range("D38").Value2 = '1.372'
$ws.Range("E38").Value2 = '  +9.93%  '

$ws.Range("D39").Value2 = '0.2306'
$ws.Range("E39").Value2 = '  -1.10%  '

$ws.Range("D40").Value2 = '0.06493'
$ws.Range("E40").Value2 = '  -1.20%  '

$ws.Range("D41").Value2 = '5.399'
$ws.Range("E41").Value2 = '  -2.94%  '

$ws.Range("D42").Value2 = '0.6330'
$ws.Range("E42").Value2 = '  -1.03%  '

$ws.Range("D43").Value2 = '11.32'
$ws.Range("E43").Value2 = '  -3.84%  '

$ws.Range("B44").Value2 = 'EnergySwap'
$ws.Range("C44").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value2 = '14.09'
$ws.Range("E44").Value2 = '  -1.74%  '

$ws.Range("B45").Value2 = 'Frax'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value2 = '1.002'
$ws.Range("E45").Value2 = '  +0.17%  '

$ws.Range("D46").Value2 = '3.791'
$ws.Range("E46").Value2 = '  +1.23%  '

$ws.Range("D47").Value2 = '0.5947'
$ws.Range("E47").Value2 = '  -0.89%  '

$ws.Range("D48").Value2 = '2.087'
$ws.Range("E48").Value2 = '  -1.92%  '

$ws.Range("D49").Value2 = '1.266'
$ws.Range("E49").Value2 = '  +0.49%  '

$ws.Range("D50").Value2 = '124.76'
$ws.Range("E50").Value2 = '  +1.24%  '

$ws.Range("D51").Value2 = '0.07320'
$ws.Range("E51").Value2 = '  +0.37%  '
